$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new survey records to the bottom of the data table
# (columns: Gender, Age, age_ed, education, Group, answers).

$ws.Range("A145").Value2 = "Мужчина"
$ws.Range("B145").Value2 = 18
$ws.Range("C145").Value2 = 12
$ws.Range("D145").Value2 = "неоконченное высшее"
$ws.Range("E145").Value2 = "A"
$ws.Range("F145").Value2 = "{4.0: ['a1', 'a1', 17044.0], 1.0: ['a2', 'a3', 12324.0], 5.0: ['a2', 'a2', 13331.0], 6.0: ['a2', 'a2', 9038.0], 2.0: ['a1', 'a3', 13251.0], 3.0: ['a2', 'a1', 41184.0]}"

$ws.Range("A146").Value2 = "Мужчина"
$ws.Range("B146").Value2 = 44
$ws.Range("C146").Value2 = 15
$ws.Range("D146").Value2 = "высшее"
$ws.Range("E146").Value2 = "B"
$ws.Range("F146").Value2 = "{4.0: ['a2', 'a1', 43169.0], 6.0: ['a2', 'a2', 33260.0], 1.0: ['a2', 'a3', 15893.0], 2.0: ['a1', 'a3', 12169.0], 5.0: ['a2', 'a2', 15221.0], 3.0: ['a1', 'a1', 12659.0]}"

# Match the author's final cursor/scroll position after the edit.
$ws.Range("F147").Select()
